$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (Yousef Haidari) - arrival/departure detection fix.
# Force text formatting first so Excel does not auto-convert the
# date/time-looking strings into date serial numbers, then restore
# the default "Normal" style so no stray number-format style is left
# attached to these cells.
$target = $ws.Range("C5:G5")
$target.NumberFormat = "@"

$ws.Range("C5").Value = "2023-07-26"
$ws.Range("D5").Value = "14:31:14"
$ws.Range("E5").Value = "2023-07-26"
$ws.Range("F5").Value = "14:31:22"
$ws.Range("G5").Value = "0:00:08"

$target.Style = "Normal"
